$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.026.42"
$ws.Range("E2").Value = "  -1.06%  "

$ws.Range("D3").Value = "2.951.84"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "377.54"
$ws.Range("E5").Value = "  +0.75%  "

$ws.Range("D6").Value = "102.17"
$ws.Range("E6").Value = "  -1.57%  "

$ws.Range("D7").Value = "0.536"
$ws.Range("E7").Value = "  -1.37%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "0.584"
$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("D10").Value = "36.38"
$ws.Range("E10").Value = "  -1.53%  "

$ws.Range("E11").Value = "  -0.35%  "

$ws.Range("D12").Value = "0.0837"
$ws.Range("E12").Value = "  -0.36%  "

$ws.Range("D13").Value = "3.414.67"
$ws.Range("E13").Value = "  +0.57%  "

$ws.Range("D14").Value = "17.88"

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.123.18"
$ws.Range("E15").Value = "  +6.46%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "7.33"
$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("D17").Value = "0.976"
$ws.Range("E17").Value = "  +4.13%  "

$ws.Range("D18").Value = "51.061.14"
$ws.Range("E18").Value = "  -0.88%  "

$ws.Range("D19").Value = "3.20"
$ws.Range("E19").Value = "  -6.76%  "

$ws.Range("D20").Value = "7.11"
$ws.Range("E20").Value = "  -3.21%  "

$ws.Range("D21").Value = "12.47"
$ws.Range("E21").Value = "  -4.13%  "

$ws.Range("D22").Value = "0.0₃0948"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").Value = "68.21"
$ws.Range("E23").Value = "  -0.34%  "

$ws.Range("D24").Value = "260.62"
$ws.Range("E24").Value = "  -0.62%  "

$ws.Range("D25").Value = "2.82"
$ws.Range("E25").Value = "  +2.03%  "

$ws.Range("D26").Value = "8.15"
$ws.Range("E26").Value = "  +10.71%  "

$ws.Range("D27").Value = "7.60"
$ws.Range("E27").Value = "  +8.18%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.167"
$ws.Range("E28").Value = "  -1.44%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("B30").Value = "LEO"
$ws.Range("C30").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D30").Value = "4.09"
$ws.Range("E30").Value = "  -1.01%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "25.57"
$ws.Range("E31").Value = "  -1.10%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.110"
$ws.Range("E32").Value = "  +9.04%  "

$ws.Range("D33").Value = "9.73"
$ws.Range("E33").Value = "  -0.81%  "

$ws.Range("E34").Value = "  -2.35%  "

$ws.Range("E35").Value = "  -2.97%  "

$ws.Range("D36").Value = "33.46"
$ws.Range("E36").Value = "  -2.22%  "

$ws.Range("D37").Value = "0.0440"
$ws.Range("E37").Value = "  +2.58%  "

$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("D39").Value = "2.94"
$ws.Range("E39").Value = "  -2.69%  "

$ws.Range("D40").Value = "16.74"
$ws.Range("E40").Value = "  -1.49%  "

$ws.Range("D41").Value = "0.114"
$ws.Range("E41").Value = "  -0.24%  "

$ws.Range("D42").Value = "2.50"
$ws.Range("E42").Value = "  -3.44%  "

$ws.Range("D43").Value = "1.76"
$ws.Range("E43").Value = "  -3.47%  "

$ws.Range("D44").Value = "121.51"
$ws.Range("E44").Value = "  -2.38%  "

$ws.Range("D45").Value = "20.88"
$ws.Range("E45").Value = "  -4.82%  "

$ws.Range("D46").Value = "2.04"
$ws.Range("E46").Value = "  -0.66%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "2.35"
$ws.Range("E47").Value = "  +1.90%  "

$ws.Range("B48").Value = "TheGraph"
$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D48").Value = "0.269"
$ws.Range("E48").Value = "  -1.78%  "

$ws.Range("D49").Value = "1.997.16"
$ws.Range("E49").Value = "  -1.42%  "

$ws.Range("D50").Value = "3.18"
$ws.Range("E50").Value = "  +0.21%  "

$ws.Range("D51").Value = "0.0329"
$ws.Range("E51").Value = "  +2.38%  "
